$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("D52").Value = 2470313.891
$ws.Range("D53").Value = 2256001.79
$ws.Range("D54").Value = 1989436.843
$ws.Range("D55").Value = 2404490.336
$ws.Range("D56").Value = 2387742.679
$ws.Range("D80").Value = 3846206.98
